# Weekly fruit/vegetable price report update ("Fruta / hortaliza, semanal")
# A new price record is inserted at row 278, shifting the existing rows
# 278:341 down to 279:342 (the row formerly at 341 becomes row 342).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 278; Excel shifts row 278 (and
# everything below it) down by one, carrying along formatting (e.g. the
# date style on column D).
$ws.Rows(278).Insert()

# Populate the newly inserted row 278 with the new weekly record.
$ws.Range("A278").Value = 10
$ws.Range("B278").Value = "Vega Modelo de Temuco"
$ws.Range("C278").Value = "La Araucanía"
$ws.Range("D278").Value = 44782
$ws.Range("E278").Value = 9
$ws.Range("F278").Value = 100112017
$ws.Range("G278").Value = "Apio"
$ws.Range("H278").Value = "Americana (o)"
$ws.Range("I278").Value = "Primera"
$ws.Range("J278").Value = 30
$ws.Range("K278").Value = 12000
$ws.Range("L278").Value = 12000
$ws.Range("M278").Value = 12000
$ws.Range("N278").Value = "$/docena de matas"
$ws.Range("O278").Value = "Provincia del Elquí"
$ws.Range("P278").Value = 2000
$ws.Range("Q278").Value = 6
$ws.Range("R278").Value = "Hortaliza"
